$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Motor Control" section header row (row 23) ---
# Merge first, then copy the formatting of the existing "AT tiny" section
# header (row 15, bold "Check Cell" style with a double border) onto it so
# all five cells end up sharing a single uniform style, then set the text.
$ws.Range("A23:E23").Merge() | Out-Null
$ws.Range("A15:E15").Copy() | Out-Null
$ws.Range("A23:E23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A23").Value = "Motor Control (Mike will have to fill this in)"

# Row heights to mirror the thick border spacing used around row 15's header
$ws.Rows.Item(22).RowHeight = 15.75
$ws.Rows.Item(23).RowHeight = 16.5
$ws.Rows.Item(24).RowHeight = 15.75

# --- New motor control pin rows (24-30) ---
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 6

$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 9

$ws.Range("C26").Value = 28
$ws.Range("D26").Value = 1

$ws.Range("C27").Value = 15
$ws.Range("D27").Value = 22

$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 27

$ws.Range("C29").Value = 11
$ws.Range("D29").Value = 17

$ws.Range("C30").Value = 12
$ws.Range("D30").Value = 18

# --- Misc view state (matches the selected cell recorded in the file) ---
$ws.Range("J17").Select() | Out-Null
